$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9208168983459473
$ws.Range("B1").Value = 2.01803731918335
$ws.Range("C1").Value = 8.919893264770508
$ws.Range("D1").Value = 1.850359201431274
$ws.Range("E1").Value = 1.426436424255371
